# Insert a new data row before the current row 218 (shifting all rows
# 218..272 down to 219..273) and populate the new row with the values
# described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 218..272 down by one, creating a blank row at 218.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218.
$ws.Range("A218").Value = 3
$ws.Range("B218").Value = "Femacal de La Calera"
$ws.Range("C218").Value = "Coquimbo"
$ws.Range("D218").Value = 44463
$ws.Range("E218").Value = 5
$ws.Range("F218").Value = 100112021
$ws.Range("G218").Value = "Ají"
$ws.Range("H218").Value = "Inferno"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 35
$ws.Range("K218").Value = 42000
$ws.Range("L218").Value = 43000
$ws.Range("M218").Value = 42429
$ws.Range("N218").Value = "`$/caja 15 kilos"
$ws.Range("O218").Value = "Región de Arica y Parinacota"
$ws.Range("P218").Value = 2829
$ws.Range("Q218").Value = 15
$ws.Range("R218").Value = "Hortaliza"
